# Auto update stock data: bump the as-of date on each block from 2025/11/15
# to 2025/11/16, and refresh Alro Steel's EBITDA figure (row 38).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $dateRows) {
    $cell = $ws.Range("A$r")
    # Force text so the yyyy/mm/dd-looking string isn't reinterpreted as a
    # date serial, then drop back to the Normal style so no stray number
    # format sticks to the cell - matches the original plain-text value.
    $cell.NumberFormat = "@"
    $cell.Value = "2025/11/16"
    $cell.Style = "Normal"
}

$b38 = $ws.Range("B38")
$b38.NumberFormat = "@"
$b38.Value = "38.59"
$b38.Style = "Normal"
